$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# Add the two new header labels in columns D and E.
$ws.Range("D1").Value = "Predicted Favor"
$ws.Range("E1").Value = "Predicted Oppose"

# Copy the bold/centered/bordered header style from an existing header cell (B1)
# onto the two new header cells so they match the look of Favor/Oppose headers.
$ws.Range("B1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# Remove the old "year" label cell in A1 entirely (content + style), since the
# year column no longer carries a header label.
$ws.Range("A1").Clear()

# --- Shift the existing "Oppose" data (old column C, rows 4:33) into column D ---
$ws.Range("D4:D33").Value2 = $ws.Range("C4:C33").Value2
$ws.Range("C4:C33").ClearContents()

# --- Append the new predicted rows (34:39) for years 2019-2024 ---
$years = @(2019, 2020, 2021, 2022, 2023, 2024)
$predictedFavor = @(0.6336356780376812, 0.6323061555557529, 0.6311282812561586, 0.6300021177793673, 0.6288935872814473, 0.6277910694889344)
$predictedOppose = @(0.3663643056059697, 0.3676938167690059, 0.3688716812173663, 0.3699978352581099, 0.3711063564324797, 0.3722088649298879)

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = 34 + $i
    $ws.Cells.Item($r, 1).Value = $years[$i]
    $ws.Cells.Item($r, 3).Value = $predictedFavor[$i]
    $ws.Cells.Item($r, 5).Value = $predictedOppose[$i]
}

# Give the new year cells (A34:A39) the same bold/centered/bordered style used
# by the other year cells in column A (e.g. A2).
$ws.Range("A2").Copy()
$ws.Range("A34:A39").PasteSpecial(-4122)
